# Applies the cell value updates described in the commit diff
# ("Atualizando o arquivo XLSX") to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("BD3").Value = 126

# Row 5
$ws.Range("G5").Value = 2.5
$ws.Range("H5").Value = 2.88
$ws.Range("I5").Value = 3.2
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("AJ5").Value = 12

# Row 7
$ws.Range("Q7").Value = 1.75
$ws.Range("U7").Value = 1.67

# Row 8
$ws.Range("U8").Value = 1.8
$ws.Range("V8").Value = 1.91

# Row 9
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91

# Row 10
$ws.Range("O10").Value = 1.14
$ws.Range("U10").Value = 2.3
$ws.Range("V10").Value = 1.59

# Row 11
$ws.Range("Q11").Value = 1.53

# Row 19
$ws.Range("G19").Value = 2.25
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 2.88
$ws.Range("W19").Value = 9.5
$ws.Range("X19").Value = 12
$ws.Range("Z19").Value = 21
$ws.Range("AI19").Value = 17
$ws.Range("AM19").Value = 29
$ws.Range("AO19").Value = 12
$ws.Range("BA19").Value = 67

# Row 22
$ws.Range("J22").Value = 1.77
$ws.Range("K22").Value = 2.87

# Row 24
$ws.Range("G24").Value = 2.6
$ws.Range("H24").Value = 3.2
$ws.Range("I24").Value = 2.7
$ws.Range("J24").Value = 3.4
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 3.5
$ws.Range("M24").Value = 1.08
$ws.Range("N24").Value = 8
$ws.Range("O24").Value = 1.4
$ws.Range("P24").Value = 2.75
$ws.Range("Q24").Value = 2.25
$ws.Range("S24").Value = 1.5
$ws.Range("T24").Value = 2.5
$ws.Range("U24").Value = 1.91
$ws.Range("V24").Value = 1.8
$ws.Range("W24").Value = 7.5
$ws.Range("X24").Value = 12
$ws.Range("Y24").Value = 10
$ws.Range("Z24").Value = 26
$ws.Range("AA24").Value = 23
$ws.Range("AB24").Value = 34
$ws.Range("AC24").Value = 8
$ws.Range("AD24").Value = 6
$ws.Range("AE24").Value = 17
$ws.Range("AF24").Value = 51
$ws.Range("AG24").Value = 351
$ws.Range("AH24").Value = 7.5
$ws.Range("AI24").Value = 12
$ws.Range("AJ24").Value = 11
$ws.Range("AK24").Value = 29
$ws.Range("AL24").Value = 23
$ws.Range("AM24").Value = 34
$ws.Range("AN24").Value = 4.5
$ws.Range("AP24").Value = 29
$ws.Range("AQ24").Value = 51
$ws.Range("AR24").Value = 81
$ws.Range("AS24").Value = 201
$ws.Range("AT24").Value = 2.5
$ws.Range("AU24").Value = 8.5
$ws.Range("AV24").Value = 67
$ws.Range("AW24").Value = 4.75
$ws.Range("AX24").Value = 17
$ws.Range("AY24").Value = 29
$ws.Range("AZ24").Value = 51
$ws.Range("BA24").Value = 81
$ws.Range("BB24").Value = 201

# Row 30
$ws.Range("M30").Value = 1.03
$ws.Range("O30").Value = 1.25

# Row 32
$ws.Range("M32").Value = 1.02
$ws.Range("O32").Value = 1.11

# Row 38
$ws.Range("N38").Value = 8

# Row 39
$ws.Range("Q39").Value = 1.95
$ws.Range("R39").Value = 1.9
